$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country names whose rank order changed (Malaui overtook Libia,
#     Republica de Chipre overtook Georgia, Santa Lucia overtook Timor Oriental) ---
$ws.Range("A108").Value = "Malaui"
$ws.Range("A109").Value = "Libia"

$ws.Range("A145").Value = "Republica de Chipre"
$ws.Range("A146").Value = "Georgia"

$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"

# --- Updated case numbers ---
# Estados Unidos
$ws.Range("B4").Value = 5007514
$ws.Range("C4").Value = 33946
$ws.Range("D4").Value = 2554225
$ws.Range("E4").Value = 2291067
$ws.Range("G4").Value = 621
$ws.Range("H4").Value = 162222

# India
$ws.Range("B6").Value = 2025409
$ws.Range("C6").Value = 62170
$ws.Range("E6").Value = 606387
$ws.Range("G6").Value = 899
$ws.Range("H6").Value = 41638

# España
$ws.Range("B12").Value = 354530
$ws.Range("C12").Value = 1683
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 28500

# row 22
$ws.Range("B22").Value = 215153
$ws.Range("C22").Value = 1049
$ws.Range("E22").Value = 9701
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = 9252

# row 23
$ws.Range("B23").Value = 195633
$ws.Range("C23").Value = 1604
$ws.Range("E23").Value = 83155
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 30312

# row 65
$ws.Range("B65").Value = 26372
$ws.Range("C65").Value = 69
$ws.Range("E65").Value = 1240
$ws.Range("G65").Value = 5
$ws.Range("H65").Value = 1768

# row 108 (now Malaui)
$ws.Range("B108").Value = 4491
$ws.Range("C108").Value = 65
$ws.Range("D108").Value = 2137
$ws.Range("E108").Value = 2217
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 137

# row 109 (now Libia)
$ws.Range("B109").Value = 4475
$ws.Range("D109").Value = 640
$ws.Range("E109").Value = 3736
$ws.Range("H109").Value = 99

# row 145 (now Republica de Chipre)
$ws.Range("B145").Value = 1208
$ws.Range("C145").Value = 13
$ws.Range("D145").Value = 856
$ws.Range("E145").Value = 333
$ws.Range("H145").Value = 19

# row 146 (now Georgia)
$ws.Range("B146").Value = 1206
$ws.Range("C146").Value = 9
$ws.Range("D146").Value = 987
$ws.Range("E146").Value = 202
$ws.Range("H146").Value = 17

# --- Updated "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Agosto de 2020 a las 21:18"
